$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 7.759382333333334
$ws.Cells.Item(2, 8).Value = 23.278147
$ws.Cells.Item(2, 9).Value = 0.03884312320086166
$ws.Cells.Item(2, 10).Value = 0.04014625174390325
$ws.Cells.Item(2, 13).Value = 121.928739
$ws.Cells.Item(2, 14).Value = 365.786217
$ws.Cells.Item(2, 15).Value = 0.2282232151508951
$ws.Cells.Item(2, 16).Value = 0.2419720431319445
$ws.Cells.Item(2, 17).Value = 946.0917033222111
$ws.Cells.Item(2, 18).Value = 8514.825329899899
$ws.Cells.Item(2, 19).Value = 0.008864902463402975
$ws.Cells.Item(2, 20).Value = 0.009714270558561661

$ws.Cells.Item(3, 7).Value = 7.759382333333334
$ws.Cells.Item(3, 8).Value = 23.278147
$ws.Cells.Item(3, 9).Value = 0.03884312320086166
$ws.Cells.Item(3, 10).Value = 0.04014625174390325
$ws.Cells.Item(3, 15).Value = 0.2768624053389947
$ws.Cells.Item(3, 16).Value = 0.2935413991166814
$ws.Cells.Item(3, 17).Value = 1147.723839048837
$ws.Cells.Item(3, 18).Value = 10329.51455143953
$ws.Cells.Item(3, 19).Value = 0.01075420052026947
$ws.Cells.Item(3, 20).Value = 0.01178458690619587

$ws.Cells.Item(4, 7).Value = 7.759382333333334
$ws.Cells.Item(4, 8).Value = 23.278147
$ws.Cells.Item(4, 9).Value = 0.03884312320086166
$ws.Cells.Item(4, 10).Value = 0.04014625174390325
$ws.Cells.Item(4, 13).Value = 83.50496933333334
$ws.Cells.Item(4, 14).Value = 250.514908
$ws.Cells.Item(4, 15).Value = 0.1563025480180701
$ws.Cells.Item(4, 16).Value = 0.1657186665504434
$ws.Cells.Item(4, 17).Value = 647.9469837906086
$ws.Cells.Item(4, 18).Value = 5831.522854115476
$ws.Cells.Item(4, 19).Value = 0.006071279129274491
$ws.Cells.Item(4, 20).Value = 0.006652983305998059

$ws.Cells.Item(5, 7).Value = 7.759382333333334
$ws.Cells.Item(5, 8).Value = 23.278147
$ws.Cells.Item(5, 9).Value = 0.03884312320086166
$ws.Cells.Item(5, 10).Value = 0.04014625174390325
$ws.Cells.Item(5, 13).Value = 91.06846250000001
$ws.Cells.Item(5, 14).Value = 182.136925
$ws.Cells.Item(5, 15).Value = 0.1704597085236707
$ws.Cells.Item(5, 16).Value = 0.1204857969594293
$ws.Cells.Item(5, 17).Value = 706.6350190463293
$ws.Cells.Item(5, 18).Value = 4239.810114277975
$ws.Cells.Item(5, 19).Value = 0.006621187458967911
$ws.Cells.Item(5, 20).Value = 0.004837053136298063

$ws.Cells.Item(6, 7).Value = 7.759382333333334
$ws.Cells.Item(6, 8).Value = 23.278147
$ws.Cells.Item(6, 9).Value = 0.03884312320086166
$ws.Cells.Item(6, 10).Value = 0.04014625174390325
$ws.Cells.Item(6, 13).Value = 89.83563
$ws.Cells.Item(6, 14).Value = 269.50689
$ws.Cells.Item(6, 15).Value = 0.1681521229683693
$ws.Cells.Item(6, 16).Value = 0.1782820942415013
$ws.Cells.Item(6, 17).Value = 697.06900032587
$ws.Cells.Item(6, 18).Value = 6273.62100293283
$ws.Cells.Item(6, 19).Value = 0.006531553628946808
$ws.Cells.Item(6, 20).Value = 0.007157357836849595

$ws.Cells.Item(7, 9).Value = 0.8631909770948131
$ws.Cells.Item(7, 10).Value = 0.8921497401307179
$ws.Cells.Item(7, 13).Value = 121.928739
$ws.Cells.Item(7, 14).Value = 365.786217
$ws.Cells.Item(7, 15).Value = 0.2282232151508951
$ws.Cells.Item(7, 16).Value = 0.2419720431319445
$ws.Cells.Item(7, 17).Value = 21024.51488231202
$ws.Cells.Item(7, 18).Value = 189220.6339408082
$ws.Cells.Item(7, 19).Value = 0.1970002200818209
$ws.Cells.Item(7, 20).Value = 0.2158752953990632

$ws.Cells.Item(8, 9).Value = 0.8631909770948131
$ws.Cells.Item(8, 10).Value = 0.8921497401307179
$ws.Cells.Item(8, 15).Value = 0.2768624053389947
$ws.Cells.Item(8, 16).Value = 0.2935413991166814
$ws.Cells.Item(8, 17).Value = 25505.28331464341
$ws.Cells.Item(8, 19).Value = 0.238985130185387
$ws.Cells.Item(8, 20).Value = 0.2618828829395547

$ws.Cells.Item(9, 9).Value = 0.8631909770948131
$ws.Cells.Item(9, 10).Value = 0.8921497401307179
$ws.Cells.Item(9, 13).Value = 83.50496933333334
$ws.Cells.Item(9, 14).Value = 250.514908
$ws.Cells.Item(9, 15).Value = 0.1563025480180701
$ws.Cells.Item(9, 16).Value = 0.1657186665504434
$ws.Cells.Item(9, 17).Value = 14398.99637193554
$ws.Cells.Item(9, 18).Value = 129590.9673474199
$ws.Cells.Item(9, 19).Value = 0.1349189491461268
$ws.Cells.Item(9, 20).Value = 0.1478458652977871

$ws.Cells.Item(10, 9).Value = 0.8631909770948131
$ws.Cells.Item(10, 10).Value = 0.8921497401307179
$ws.Cells.Item(10, 13).Value = 91.06846250000001
$ws.Cells.Item(10, 14).Value = 182.136925
$ws.Cells.Item(10, 15).Value = 0.1704597085236707
$ws.Cells.Item(10, 16).Value = 0.1204857969594293
$ws.Cells.Item(10, 17).Value = 15703.19073947386
$ws.Cells.Item(10, 18).Value = 94219.14443684313
$ws.Cells.Item(10, 19).Value = 0.1471392823558444
$ws.Cells.Item(10, 20).Value = 0.1074913724467973

$ws.Cells.Item(11, 9).Value = 0.8631909770948131
$ws.Cells.Item(11, 10).Value = 0.8921497401307179
$ws.Cells.Item(11, 13).Value = 89.83563
$ws.Cells.Item(11, 14).Value = 269.50689
$ws.Cells.Item(11, 15).Value = 0.1681521229683693
$ws.Cells.Item(11, 16).Value = 0.1782820942415013
$ws.Cells.Item(11, 17).Value = 15490.60996929425
$ws.Cells.Item(11, 18).Value = 139415.4897236482
$ws.Cells.Item(11, 19).Value = 0.1451473953256338
$ws.Cells.Item(11, 20).Value = 0.1590543240475155

$ws.Cells.Item(12, 7).Value = 0.05240566666666666
$ws.Cells.Item(12, 8).Value = 0.157217
$ws.Cells.Item(12, 9).Value = 0.0002623404388789996
$ws.Cells.Item(12, 10).Value = 0.0002711415672571033
$ws.Cells.Item(12, 13).Value = 121.928739
$ws.Cells.Item(12, 14).Value = 365.786217
$ws.Cells.Item(12, 15).Value = 0.2282232151508951
$ws.Cells.Item(12, 16).Value = 0.2419720431319445
$ws.Cells.Item(12, 17).Value = 6.389756853120999
$ws.Cells.Item(12, 18).Value = 57.50781167808899
$ws.Cells.Item(12, 19).Value = 0.00005987217842506216
$ws.Cells.Item(12, 20).Value = 0.00006560867900719884

$ws.Cells.Item(13, 7).Value = 0.05240566666666666
$ws.Cells.Item(13, 8).Value = 0.157217
$ws.Cells.Item(13, 9).Value = 0.0002623404388789996
$ws.Cells.Item(13, 10).Value = 0.0002711415672571033
$ws.Cells.Item(13, 15).Value = 0.2768624053389947
$ws.Cells.Item(13, 16).Value = 0.2935413991166814
$ws.Cells.Item(13, 17).Value = 7.751549073203331
$ws.Cells.Item(13, 18).Value = 69.76394165882999
$ws.Cells.Item(13, 19).Value = 0.00007263220492572734
$ws.Cells.Item(13, 20).Value = 0.00007959127501133988

$ws.Cells.Item(14, 7).Value = 0.05240566666666666
$ws.Cells.Item(14, 8).Value = 0.157217
$ws.Cells.Item(14, 9).Value = 0.0002623404388789996
$ws.Cells.Item(14, 10).Value = 0.0002711415672571033
$ws.Cells.Item(14, 13).Value = 83.50496933333334
$ws.Cells.Item(14, 14).Value = 250.514908
$ws.Cells.Item(14, 15).Value = 0.1563025480180701
$ws.Cells.Item(14, 16).Value = 0.1657186665504434
$ws.Cells.Item(14, 17).Value = 4.376133587892888
$ws.Cells.Item(14, 18).Value = 39.38520229103599
$ws.Cells.Item(14, 19).Value = 0.0000410044790449664
$ws.Cells.Item(14, 20).Value = 0.00004493321897224452

$ws.Cells.Item(15, 7).Value = 0.05240566666666666
$ws.Cells.Item(15, 8).Value = 0.157217
$ws.Cells.Item(15, 9).Value = 0.0002623404388789996
$ws.Cells.Item(15, 10).Value = 0.0002711415672571033
$ws.Cells.Item(15, 13).Value = 91.06846250000001
$ws.Cells.Item(15, 14).Value = 182.136925
$ws.Cells.Item(15, 15).Value = 0.1704597085236707
$ws.Cells.Item(15, 16).Value = 0.1204857969594293
$ws.Cells.Item(15, 17).Value = 4.772503489620833
$ws.Cells.Item(15, 18).Value = 28.635020937725
$ws.Cells.Item(15, 19).Value = 0.00004471847474528612
$ws.Cells.Item(15, 20).Value = 0.0000326687078198008

$ws.Cells.Item(16, 7).Value = 0.05240566666666666
$ws.Cells.Item(16, 8).Value = 0.157217
$ws.Cells.Item(16, 9).Value = 0.0002623404388789996
$ws.Cells.Item(16, 10).Value = 0.0002711415672571033
$ws.Cells.Item(16, 13).Value = 89.83563
$ws.Cells.Item(16, 14).Value = 269.50689
$ws.Cells.Item(16, 15).Value = 0.1681521229683693
$ws.Cells.Item(16, 16).Value = 0.1782820942415013
$ws.Cells.Item(16, 17).Value = 4.707896080569999
$ws.Cells.Item(16, 18).Value = 42.37106472513
$ws.Cells.Item(16, 19).Value = 0.0000441131017379575
$ws.Cells.Item(16, 20).Value = 0.00004833968644651925

$ws.Cells.Item(17, 7).Value = 19.452549
$ws.Cells.Item(17, 8).Value = 38.905098
$ws.Cells.Item(17, 9).Value = 0.09737859599105524
$ws.Cells.Item(17, 10).Value = 0.06709700125311635
$ws.Cells.Item(17, 13).Value = 121.928739
$ws.Cells.Item(17, 14).Value = 365.786217
$ws.Cells.Item(17, 15).Value = 0.2282232151508951
$ws.Cells.Item(17, 16).Value = 0.2419720431319445
$ws.Cells.Item(17, 17).Value = 2371.824769905711
$ws.Cells.Item(17, 18).Value = 14230.94861943426
$ws.Cells.Item(17, 19).Value = 0.02222405626395869
$ws.Cells.Item(17, 20).Value = 0.01623559848124321

$ws.Cells.Item(18, 7).Value = 19.452549
$ws.Cells.Item(18, 8).Value = 38.905098
$ws.Cells.Item(18, 9).Value = 0.09737859599105524
$ws.Cells.Item(18, 10).Value = 0.06709700125311635
$ws.Cells.Item(18, 15).Value = 0.2768624053389947
$ws.Cells.Item(18, 16).Value = 0.2935413991166814
$ws.Cells.Item(18, 17).Value = 2877.310752127169
$ws.Cells.Item(18, 18).Value = 17263.86451276302
$ws.Cells.Item(18, 19).Value = 0.02696047231461774
$ws.Cells.Item(18, 20).Value = 0.0196957476243735

$ws.Cells.Item(19, 7).Value = 19.452549
$ws.Cells.Item(19, 8).Value = 38.905098
$ws.Cells.Item(19, 9).Value = 0.09737859599105524
$ws.Cells.Item(19, 10).Value = 0.06709700125311635
$ws.Cells.Item(19, 13).Value = 83.50496933333334
$ws.Cells.Item(19, 14).Value = 250.514908
$ws.Cells.Item(19, 15).Value = 0.1563025480180701
$ws.Cells.Item(19, 16).Value = 0.1657186665504434
$ws.Cells.Item(19, 17).Value = 1624.384507700164
$ws.Cells.Item(19, 18).Value = 9746.307046200982
$ws.Cells.Item(19, 19).Value = 0.01522052267582416
$ws.Cells.Item(19, 20).Value = 0.01111922557719987

$ws.Cells.Item(20, 7).Value = 19.452549
$ws.Cells.Item(20, 8).Value = 38.905098
$ws.Cells.Item(20, 9).Value = 0.09737859599105524
$ws.Cells.Item(20, 10).Value = 0.06709700125311635
$ws.Cells.Item(20, 13).Value = 91.06846250000001
$ws.Cells.Item(20, 14).Value = 182.136925
$ws.Cells.Item(20, 15).Value = 0.1704597085236707
$ws.Cells.Item(20, 16).Value = 0.1204857969594293
$ws.Cells.Item(20, 17).Value = 1771.513729135912
$ws.Cells.Item(20, 18).Value = 7086.05491654365
$ws.Cells.Item(20, 19).Value = 0.01659912708907957
$ws.Cells.Item(20, 20).Value = 0.008084235669569553

$ws.Cells.Item(21, 7).Value = 19.452549
$ws.Cells.Item(21, 8).Value = 38.905098
$ws.Cells.Item(21, 9).Value = 0.09737859599105524
$ws.Cells.Item(21, 10).Value = 0.06709700125311635
$ws.Cells.Item(21, 13).Value = 89.83563
$ws.Cells.Item(21, 14).Value = 269.50689
$ws.Cells.Item(21, 15).Value = 0.1681521229683693
$ws.Cells.Item(21, 16).Value = 0.1782820942415013
$ws.Cells.Item(21, 17).Value = 1747.53199452087
$ws.Cells.Item(21, 18).Value = 10485.19196712522
$ws.Cells.Item(21, 19).Value = 0.01637441764757507
$ws.Cells.Item(21, 20).Value = 0.01196219390073022

$ws.Cells.Item(22, 7).Value = 0.06491533333333334
$ws.Cells.Item(22, 8).Value = 0.194746
$ws.Cells.Item(22, 9).Value = 0.0003249632743909987
$ws.Cells.Item(22, 10).Value = 0.0003358653050055137
$ws.Cells.Item(22, 13).Value = 121.928739
$ws.Cells.Item(22, 14).Value = 365.786217
$ws.Cells.Item(22, 15).Value = 0.2282232151508951
$ws.Cells.Item(22, 16).Value = 0.2419720431319445
$ws.Cells.Item(22, 17).Value = 7.915044735098
$ws.Cells.Item(22, 18).Value = 71.235402615882
$ws.Cells.Item(22, 19).Value = 0.00007416416328747626
$ws.Cells.Item(22, 20).Value = 0.00008127001406931785

$ws.Cells.Item(23, 7).Value = 0.06491533333333334
$ws.Cells.Item(23, 8).Value = 0.194746
$ws.Cells.Item(23, 9).Value = 0.0003249632743909987
$ws.Cells.Item(23, 10).Value = 0.0003358653050055137
$ws.Cells.Item(23, 15).Value = 0.2768624053389947
$ws.Cells.Item(23, 16).Value = 0.2935413991166814
$ws.Cells.Item(23, 17).Value = 9.601908036726666
$ws.Cells.Item(23, 18).Value = 86.41717233054
$ws.Cells.Item(23, 19).Value = 0.00008997011379472765
$ws.Cells.Item(23, 20).Value = 0.00009859037154606943

$ws.Cells.Item(24, 7).Value = 0.06491533333333334
$ws.Cells.Item(24, 8).Value = 0.194746
$ws.Cells.Item(24, 9).Value = 0.0003249632743909987
$ws.Cells.Item(24, 10).Value = 0.0003358653050055137
$ws.Cells.Item(24, 13).Value = 83.50496933333334
$ws.Cells.Item(24, 14).Value = 250.514908
$ws.Cells.Item(24, 15).Value = 0.1563025480180701
$ws.Cells.Item(24, 16).Value = 0.1657186665504434
$ws.Cells.Item(24, 17).Value = 5.420752919263112
$ws.Cells.Item(24, 18).Value = 48.786776273368
$ws.Cells.Item(24, 19).Value = 0.00005079258779960836
$ws.Cells.Item(24, 20).Value = 0.00005565915048607167

$ws.Cells.Item(25, 7).Value = 0.06491533333333334
$ws.Cells.Item(25, 8).Value = 0.194746
$ws.Cells.Item(25, 9).Value = 0.0003249632743909987
$ws.Cells.Item(25, 10).Value = 0.0003358653050055137
$ws.Cells.Item(25, 13).Value = 91.06846250000001
$ws.Cells.Item(25, 14).Value = 182.136925
$ws.Cells.Item(25, 15).Value = 0.1704597085236707
$ws.Cells.Item(25, 16).Value = 0.1204857969594293
$ws.Cells.Item(25, 17).Value = 5.911739599341668
$ws.Cells.Item(25, 18).Value = 35.47043759605
$ws.Cells.Item(25, 19).Value = 0.00005539314503358728
$ws.Cells.Item(25, 20).Value = 0.00004046699894461112

$ws.Cells.Item(26, 7).Value = 0.06491533333333334
$ws.Cells.Item(26, 8).Value = 0.194746
$ws.Cells.Item(26, 9).Value = 0.0003249632743909987
$ws.Cells.Item(26, 10).Value = 0.0003358653050055137
$ws.Cells.Item(26, 13).Value = 89.83563
$ws.Cells.Item(26, 14).Value = 269.50689
$ws.Cells.Item(26, 15).Value = 0.1681521229683693
$ws.Cells.Item(26, 16).Value = 0.1782820942415013
$ws.Cells.Item(26, 17).Value = 5.83170986666
$ws.Cells.Item(26, 18).Value = 52.48538879994
$ws.Cells.Item(26, 19).Value = 0.00005464326447559915
$ws.Cells.Item(26, 20).Value = 0.00005987876995944356
